$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text (string) representation
# instead of being auto-converted to numbers by Excel when values look numeric
# (e.g. "1.000", "30.438.08", "12.97").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.438.08'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '1.916.42'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '243.69'
$ws.Range('E5').Value = '  +1.92%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '0.4694'
$ws.Range('E7').Value = '  -1.30%  '
$ws.Range('D8').Value = '0.2864'
$ws.Range('E8').Value = '  -1.09%  '
$ws.Range('D9').Value = '0.06837'
$ws.Range('E9').Value = '  +3.86%  '
$ws.Range('D10').Value = '110.38'
$ws.Range('E10').Value = '  +11.39%  '
$ws.Range('D11').Value = '18.40'
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').Value = '0.07739'
$ws.Range('E12').Value = '  +1.84%  '
$ws.Range('D13').Value = '1.894.18'
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('D15').Value = '0.6570'
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').Value = '295.69'
$ws.Range('E16').Value = '  -4.15%  '
$ws.Range('D17').Value = '30.444.25'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '0.000007641'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = '12.97'
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '0.9995'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').Value = '2.142.42'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '5.246'
$ws.Range('E23').Value = '  +2.34%  '
$ws.Range('D24').Value = '6.214'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = '9.367'
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '21.72'
$ws.Range('E26').Value = '  +5.54%  '
$ws.Range('D27').Value = '168.87'
$ws.Range('E27').Value = '  +0.85%  '
$ws.Range('D28').Value = '2.093'
$ws.Range('E28').Value = '  +7.37%  '
$ws.Range('D29').Value = '0.1068'
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('D31').Value = '4.176'
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('D32').Value = '3.984'
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').Value = '0.05044'
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '1.154'
$ws.Range('E34').Value = '  -1.49%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '0.7355'
$ws.Range('E35').Value = '  +0.79%  '
$ws.Range('D36').Value = '0.02069'
$ws.Range('E36').Value = '  +5.67%  '
$ws.Range('D37').Value = '2.741'
$ws.Range('E37').Value = '  +0.90%  '
$ws.Range('D38').Value = '2.684'
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('D39').Value = '2.058'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('D40').Value = '109.52'
$ws.Range('E40').Value = '  +1.19%  '
$ws.Range('D41').Value = '0.8704'
$ws.Range('E41').Value = '  -4.20%  '
$ws.Range('D42').Value = '5.844'
$ws.Range('E42').Value = '  +3.79%  '
$ws.Range('D43').Value = '0.4252'
$ws.Range('E43').Value = '  +0.98%  '
$ws.Range('D44').Value = '0.9997'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '51.36'
$ws.Range('E45').Value = '  +20.12%  '
$ws.Range('D46').Value = '67.41'
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('D47').Value = '7.193'
$ws.Range('E47').Value = '  -2.73%  '
$ws.Range('D48').Value = '9.227'
$ws.Range('E48').Value = '  +2.33%  '
$ws.Range('D49').Value = '0.1218'
$ws.Range('E49').Value = '  -0.95%  '
$ws.Range('D50').Value = '34.98'
$ws.Range('E50').Value = '  +0.40%  '
$ws.Range('D51').Value = '0.2440'
$ws.Range('E51').Value = '  +10.04%  '
